# Updated Formal Requirements Excel Sheet
# - Removed the "LED indicator" electrical requirement (old CHG-EE-0005 row)
#   and renumbered the following row's ID back to CHG-EE-0005.
# - Removed the "USB-C port and LED same edge" mechanical requirement
#   (old CHG-ME-0004, the last row of the table).
# - Refreshed the status line date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete bottom-most row first so the earlier row's index doesn't shift
# before we get to it.
$ws.Rows("17:17").Delete()
$ws.Rows("12:12").Delete()

# After deleting old row 12 (the LED requirement), the row that used to be
# CHG-EE-0006 slides up into row 12; relabel it back to CHG-EE-0005 so the
# ID sequence stays contiguous.
$ws.Range("A12").Value = "CHG-EE-0005"

# Update the status line with the new date.
$ws.Range("A4").Value = "Status: In Progress as of 4 September 2025"

# Re-apply the bottom-border "end of table" formatting to the new last
# data row (now row 15, since two rows were removed from the 17-row table).
$ws.Range("A15").Borders.Item(7).LineStyle = 0
$ws.Range("A15:D15").Borders.Item(9).LineStyle = 1

# Leave a blank spacer row below (row 16, no border) matching the table's
# prior trailing blank row.
$ws.Range("A16:D16").Borders.Item(9).LineStyle = 0

# Restore the cursor/selection to just below the (now shorter) table.
$ws.Range("B18").Select()
